# Sandadi_LabExam03Grading.xlsx — re-grading of the CustomerMappingDriver
# Class section (rows 65-80 in the Java driver file under review).
#
# The submission was re-evaluated: it throws an UnsupportedOperationException
# when the driver runs, so the "successfully scanning data" score moves from
# 8/16 to 9/16 and the grading-comment cells for the Driver-class rows and the
# generic "Compilation errors" row are updated to describe the new findings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 — "For successfully scanning data from input file" (CustomerMappingDriver Class)
# Actual points 8 -> 9, and the grading comment is replaced with the specific deduction note.
$ws.Range("E29").Value = 9
$ws.Range("F29").Value = '(-1) for nor declaring Customer object, (-2) for comparing the string with "Product" not "Customer", (-2) for not initializing the customer object, (-2) for not initializing product object'

# Row 37 — "Compilation errors if any" (Generic section)
$ws.Range("F37").Value = 'For getting UnsupportedOperationException when running the driver file'

# Row 30 — "For correct and properly aligned output" (CustomerMappingDriver Class)
$ws.Range("F30").Value = '(-4) for no output due to UnsupportedOperationException exception'

# Leave the cursor/selection where the grader ended up looking last: F30.
$ws.Range("F30").Select() | Out-Null
